$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 currently holds "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)"
# Row 24 currently holds "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)"
# The edit swaps their order, so the weaker-requirement course (LOT2028) now appears first (row 23)
# and the "Indicação de Conjunto" course (LOT2038) appears second (row 24).

$textConjunto = "LOT2038 -  Tecnologia de Bebidas  (Indicação de Conjunto)`n"
$textFraco = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

$ws.Range("B23").Value = $textFraco
$ws.Range("C23").Value = $textFraco

$ws.Range("B24").Value = $textConjunto
$ws.Range("C24").Value = $textConjunto
